$wb = $excel.ActiveWorkbook

# Rename headers on existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" sheet as the last sheet in the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used by the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $wsForecast.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: ds (date serial), PO_Forecast, yhat_lower, yhat_upper
$data = @(
  @(45074.99999999999, 71, 2.394083268828834, 137.3584843487138),
  @(45088.99999999999, 71, 5.579220191794879, 137.5540051080016),
  @(45095.99999999999, 70, 1.323532121150619, 133.4184869762227),
  @(45102.99999999999, 70, 6.325735021379826, 133.5615776567447),
  @(45109.99999999999, 70, 4.112812300605833, 136.3005551660014),
  @(45116.99999999999, 70, 5.655908364072148, 131.2143770495787),
  @(45263.99999999999, 67, 2.048286003629377, 133.1119854170524),
  @(45277.99999999999, 66, 0.05116035914404427, 134.7870525981313),
  @(45298.99999999999, 66, 3.512178767627892, 136.6111250726028),
  @(45312.99999999999, 66, 0.3960586357298023, 131.0568320732737),
  @(45326.99999999999, 65, 1.095424265663088, 134.2514203360018),
  @(45333.99999999999, 65, 6.410480254402464, 131.3851619170956),
  @(45340.99999999999, 65, 2.954456017892618, 131.7031150292175),
  @(45347.99999999999, 65, -2.621191052414386, 129.3340317639649),
  @(45354.99999999999, 65, -1.485581562797301, 126.7037212385848),
  @(45361.99999999999, 65, -0.1324027223177635, 128.7799430288639),
  @(45368.99999999999, 65, -2.307039606626153, 129.1681142707251),
  @(45375.99999999999, 64, -0.6109265303804738, 128.7506079544613),
  @(45382.99999999999, 64, -2.517151330767197, 122.0209512113947),
  @(45389.99999999999, 64, -1.568365583703413, 130.5921894900399),
  @(45424.99999999999, 63, 0.69291928220054, 128.3322815210329),
  @(45431.99999999999, 63, -1.458553376453992, 127.1609324403595),
  @(45445.99999999999, 63, -3.320568760931742, 128.8440122335631),
  @(45459.99999999999, 63, -2.416404718735578, 128.2946029772382),
  @(45466.99999999999, 62, -4.493820041253024, 126.922277068105),
  @(45473.99999999999, 62, -8.418140979424011, 126.7751338482055),
  @(45480.99999999999, 62, -0.1433831920465645, 123.1407905035524),
  @(45487.99999999999, 62, -5.334105357355641, 124.2272931725005),
  @(45508.99999999999, 61, -4.246188613146841, 129.2949762948291),
  @(45515.99999999999, 61, -5.242095224747846, 127.4990614708564),
  @(45543.99999999999, 61, -4.267386988017133, 125.3524820527388),
  @(45550.99999999999, 61, -2.924387409183655, 127.8372586730408),
  @(45564.99999999999, 60, -11.99562975107325, 125.3738176635839),
  @(45571.99999999999, 60, -9.087768371728973, 123.8993424016081),
  @(45585.99999999999, 60, -8.532440226993318, 120.4384538003725),
  @(45620.99999999999, 59, -8.031052848142775, 119.5280446221731),
  @(45641.99999999999, 59, -8.841620622091972, 121.7423975116792),
  @(45648.99999999999, 58, -5.622966964846689, 122.2844336782388),
  @(45655.99999999999, 58, -5.966763747005868, 123.7363605495982),
  @(45662.99999999999, 58, -7.163420104533773, 119.7489654830106),
  @(45669.99999999999, 58, -5.521460138583962, 119.9515847579159),
  @(45676.99999999999, 58, -5.773861640769336, 125.5575111177249),
  @(45683.99999999999, 58, -5.050634307452874, 121.7092409083438),
  @(45690.99999999999, 57, -11.30201849006116, 123.3503977711724),
  @(45697.99999999999, 57, -4.117597907741912, 125.1044180384825)
)

$r = 2
foreach ($row in $data) {
    $dateCell = $wsForecast.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]

    $r++
}

$wsForecast.Range("A1").Select()
